$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the current row 720, shifting the existing
# rows 720-745 down to 723-748 (this matches the diff, which shows the old
# row-720 data reappearing - unchanged - at row 723, etc.).
$ws.Rows.Item(720).Insert()
$ws.Rows.Item(720).Insert()
$ws.Rows.Item(720).Insert()

# Make sure the new rows use the same date format as the rest of column D.
$ws.Range("D720:D722").NumberFormat = "YYYY-MM-DD HH:MM:SS"

function Set-GrapeRow {
    param(
        [int]$row,
        [double]$fecha,
        [string]$variedad,
        [string]$calidad,
        [double]$volumen,
        [double]$precioMin,
        [double]$precioMax,
        [double]$precioProm,
        [string]$unidad,
        [string]$origen,
        [double]$precioKg,
        [double]$kgUnidad
    )

    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100109
    $ws.Cells.Item($row, 8).Value = "Uva"
    $ws.Cells.Item($row, 9).Value = 100109001
    $ws.Cells.Item($row, 10).Value = "Uva"
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-GrapeRow 720 45075 "Autumn Royal"      "Primera" 250 12000 12000 12000 "`$/bandeja 18 kilos" "Provincia de Limarí" 667 18
Set-GrapeRow 721 45075 "Crimpson Seedless" "Primera" 200 12000 12000 12000 "`$/bandeja 18 kilos" "Provincia de Limarí" 667 18
Set-GrapeRow 722 45075 "Red Globe"         "Primera" 180 12000 12000 12000 "`$/bandeja 18 kilos" "Provincia de Limarí" 667 18
